# Add a new "Charging stations list management" status row to the
# "Charging Stations Features" sheet (sheet 2), recording it as 1555
# with status "Processing" (mirrors the existing rows' layout).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)
$ws.Activate()

# New row 17: Sno / Task / Status (columns A, B, D - column C is unused
# on this sheet, matching every other row).
$ws.Cells.Item(17, 1).Value = 1555
$ws.Cells.Item(17, 2).Value = "Charging stations list management"

# Copy the formatting (fill colour / style) used by the other Status
# cells in column D (e.g. D16, which shows "Pending") onto the new
# D17 cell, then set its text to "Processing".
$ws.Range("D16").Copy()
$ws.Range("D17").PasteSpecial(-4122)
$ws.Range("D17").Value = "Processing"

# Match the author's final selection/active cell.
[void]$ws.Range("D17").Select()
